# Add new "DHSC" department rows to the SalaryData table (Table1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SalaryData")
$table = $ws.ListObjects.Item("Table1")

$dept = "DHSC"
$rows = @(
    @{ Grade = "AO";          Salary = 19340 },
    @{ Grade = "EO";          Salary = 23690 },
    @{ Grade = "HEO";         Salary = 28966 },
    @{ Grade = "Fast Stream"; Salary = 29521 },
    @{ Grade = "SEO";         Salary = 36819 },
    @{ Grade = "G7";          Salary = 49529 },
    @{ Grade = "G6";          Salary = 62404 }
)

$firstNewRow = $table.ListRows.Count + 1 + 1   # header occupies row 1

foreach ($row in $rows) {
    $listRow = $table.ListRows.Add()
    $r = $listRow.Range

    $r.Cells.Item(1, 1).Value = $dept
    $r.Cells.Item(1, 2).Value = $row.Grade
    $r.Cells.Item(1, 3).Value = $row.Salary
    $r.Cells.Item(1, 4).Formula = "=Table1[[#This Row], [Salary]]/52"
    $r.Cells.Item(1, 5).Formula = "=Table1[[#This Row], [Weekly]]/5"
    $r.Cells.Item(1, 6).Formula = "=Table1[[#This Row], [Weekly]]/37"
    $r.Cells.Item(1, 7).Formula = "=Table1[[#This Row], [Hr]]/60"
    $r.Cells.Item(1, 8).Formula = "=Table1[[#This Row], [Min]]/60"
}

$lastNewRow = $table.ListRows.Count + 1

# The new rows inherit the default (header-like) style; re-apply the
# established body-row formatting by copying it down from the row
# immediately above the freshly added block.
$ws.Range("C" + ($firstNewRow - 1) + ":H" + ($firstNewRow - 1)).Copy()
$ws.Range("C" + $firstNewRow + ":H" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D" + ($firstNewRow - 1) + ":H" + $lastNewRow).Select()
